# Bug fix: remove ohlc data carryover
# Clears the OPEN/HIGH/LOW/CLOSE (columns B:E) values that were
# incorrectly carried over on specific rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5,10,15,20,25,30,35,40,45,50,55,60,65,70,75,80,85,90,95,100,105,110,115,120,125,130,135,140,145,150,155,160,165,170,175,180,185,190,195,200,205,210,215,220,225,230,235,240,245,250,255,256,259,260)

foreach ($r in $rows) {
    $range = $ws.Range("B$r`:E$r")
    $range.ClearContents()
}
